# "correção de grallha no nome" - fix the typo in the author credit on the
# title slide's subtitle placeholder:
#   - "Elaborado por:" gets split into two runs: "Elaborado" + " por:"
#   - the student number is corrected: "Diogo Guimarães nº29528"
#     -> "Diogo Guimarães nº39528", split into three runs:
#     "Diogo " + "Guimarães" + " nº39528"
#
# This text lives on Slide 1, in the "Subtítulo 2" subtitle placeholder.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Locate the subtitle placeholder by name, falling back to the known
# shape index if the name lookup is unavailable for some reason.
$shp = $null
try {
    $shp = $s.Shapes.Item("Subtítulo 2")
} catch {
    $shp = $null
}
if ($shp -eq $null) {
    $shp = $s.Shapes.Item(4)
}

$tr = $shp.TextFrame.TextRange

$paraCount = $tr.Paragraphs().Count
for ($i = 1; $i -le $paraCount; $i++) {
    $para = $tr.Paragraphs($i, 1)
    # TextRange.Text includes a trailing paragraph-mark character (CR);
    # strip it before comparing/measuring the visible text.
    $ptext = $para.Text.TrimEnd([char]13)
    $pstart = $para.Start

    if ($ptext -eq "Elaborado por:") {
        # Split "Elaborado por:" -> "Elaborado" + " por:". Re-assigning a
        # sub-range's Text (even to its current value) forces the host to
        # break the run at that boundary, the same way PowerPoint itself
        # splits a run after an in-place text edit.
        $cut = "Elaborado".Length
        $tr.Characters($pstart, $cut).Text = $ptext.Substring(0, $cut)
        $tr.Characters($pstart + $cut, $ptext.Length - $cut).Text = $ptext.Substring($cut)
    }
    elseif ($ptext.StartsWith("Diogo Guimar") -and $ptext.Contains("29528")) {
        # Fix the transposed digit in the student number, then split into
        # "Diogo " + "Guimarães" + " nº39528".
        $fixed = $ptext -replace "29528", "39528"
        $cut1 = "Diogo ".Length
        $cut2 = "Guimarães".Length
        $cut3 = $fixed.Length - $cut1 - $cut2

        $tr.Characters($pstart, $cut1).Text = $fixed.Substring(0, $cut1)
        $tr.Characters($pstart + $cut1, $cut2).Text = $fixed.Substring($cut1, $cut2)
        $tr.Characters($pstart + $cut1 + $cut2, $cut3).Text = $fixed.Substring($cut1 + $cut2)
    }
}
